$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize capitalization of "de" -> "De" in a few specific cells
$ws.Range("B6").Value = "Mazapa De Madero"
$ws.Range("A9").Value = "Ciudad De México"
$ws.Range("A14").Value = "Estado De México"
$ws.Range("B17").Value = "Santa Cruz De Juventino Rosas"

# Remove the trailing metadata/footer rows (37-41), which sat below a
# blank row 36 beneath the data table that ends at row 35.
$ws.Rows("37:41").Delete()
